# Apply the Sun Nov 26 05:36:25 UTC 2023 cryptos-list refresh (GitHub Actions run).
# Price (column D) and 1h-volume-change (column E) cells are plain text in this
# sheet (coinranking.com scrape), so any D value that LOOKS like a plain number
# needs to be forced back to Text - otherwise Excel's COM layer auto-converts it
# to a Number (dropping e.g. a trailing '0') when .Value is assigned.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.813.40'
$ws.Range("E2").Value = '  +0.03%  '

$ws.Range("D3").Value = '2.082.28'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("E4").Value = '  +0.07%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '233.63'
$cell.ClearFormats()
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("E6").Value = '  +0.03%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '58.67'
$cell.ClearFormats()
$ws.Range("E7").Value = '  -0.58%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  +0.64%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0789'
$cell.ClearFormats()
$ws.Range("E10").Value = '  +0.02%  '

$ws.Range("E11").Value = '  +3.45%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '14.99'
$cell.ClearFormats()
$ws.Range("E12").Value = '  +2.00%  '

$ws.Range("D13").Value = '2.390.04'
$ws.Range("E13").Value = '  -0.02%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '21.34'
$cell.ClearFormats()
$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("E15").Value = '  +1.11%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '5.38'
$cell.ClearFormats()
$ws.Range("E16").Value = '  +1.69%  '

$ws.Range("D17").Value = '2.089.93'
$ws.Range("E17").Value = '  +0.35%  '

$ws.Range("D18").Value = '37.669.68'
$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("E19").Value = '  -0.44%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '71.53'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +0.16%  '

$ws.Range("D21").Value = '0.0₃0843'
$ws.Range("E21").Value = '  +1.53%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '230.31'
$cell.ClearFormats()
$ws.Range("E22").Value = '  +0.60%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.39'
$cell.ClearFormats()
$ws.Range("E24").Value = '  -0.72%  '

$ws.Range("E25").Value = '  +1.43%  '

$ws.Range("E26").Value = '  +11.05%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '172.06'
$cell.ClearFormats()
$ws.Range("E27").Value = '  +1.05%  '

$ws.Range("E28").Value = '  -1.76%  '

$ws.Range("E29").Value = '  +0.12%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '19.54'
$cell.ClearFormats()
$ws.Range("E30").Value = '  -0.09%  '

$ws.Range("E31").Value = '  +1.30%  '

$ws.Range("E32").Value = '  +0.73%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.0634'
$cell.ClearFormats()
$ws.Range("E33").Value = '  +0.66%  '

$ws.Range("E34").Value = '  -1.15%  '

$ws.Range("E35").Value = '  -1.75%  '

$ws.Range("E36").Value = '  -0.71%  '

$ws.Range("E37").Value = '  -1.33%  '

$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("E39").Value = '  +0.57%  '

$ws.Range("E40").Value = '  +9.24%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '101.64'
$cell.ClearFormats()
$ws.Range("E41").Value = '  +2.94%  '

$ws.Range("E42").Value = '  -1.64%  '

$ws.Range("E43").Value = '  -0.56%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '16.90'
$cell.ClearFormats()
$ws.Range("E44").Value = '  +5.21%  '

$ws.Range("D45").Value = '1.451.72'
$ws.Range("E45").Value = '  -0.48%  '

$ws.Range("E46").Value = '  -0.46%  '

$ws.Range("E47").Value = '  -0.15%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '4.10'
$cell.ClearFormats()
$ws.Range("E48").Value = '  -5.60%  '

$ws.Range("E49").Value = '  -0.72%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.99'
$cell.ClearFormats()
$ws.Range("E50").Value = '  -1.41%  '

$ws.Range("D51").Value = '2.274.95'
$ws.Range("E51").Value = '  -0.02%  '
